$wb = $excel.ActiveWorkbook

# --- Text update: "Ready for handoff" -> "In Translation" -------------
# This status string appears on all three sheets: the "Overview" summary
# sheet (columns E/F, one per locale) and on each per-locale detail sheet
# (column C, "Status").
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Column width update ------------------------------------------------
# Narrow the status columns on all three sheets to match the new
# (shorter) status text column width.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
